$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates: recomputed TPM-derived specificity values, new row 3 inserted below ---
$ws.Cells.Item(2, 8).Value = 0.66412
$ws.Cells.Item(2, 9).Value = 0.9279112495895713
$ws.Cells.Item(2, 10).Value = 0.9279112495895712
$ws.Cells.Item(2, 19).Value = 0.9279112495895713
$ws.Cells.Item(2, 20).Value = 0.9279112495895712

# --- Row 3: new Resolving-Mac -> Pnoc/Oprl1 -> FAPs edge ---
$ws.Cells.Item(3, 1).Value = "Resolving-Mac"
$ws.Cells.Item(3, 2).Value = "Pnoc"
$ws.Cells.Item(3, 3).Value = "Oprl1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01719833333333333
$ws.Cells.Item(3, 8).Value = 0.051595
$ws.Cells.Item(3, 9).Value = 0.07208875041042873
$ws.Cells.Item(3, 10).Value = 0.07208875041042873
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.7589613333333333
$ws.Cells.Item(3, 14).Value = 2.276884
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.01305286999777778
$ws.Cells.Item(3, 18).Value = 0.11747582998
$ws.Cells.Item(3, 19).Value = 0.07208875041042873
$ws.Cells.Item(3, 20).Value = 0.07208875041042873
